$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 142146.52
$ws.Range("I15").Value = 142146.52
$ws.Range("K15").Value = 426439.5599999999
$ws.Range("M15").Value = -426270.5599999999

$ws.Range("H69").Value = 3839.6592
$ws.Range("J69").Value = 6163.5713
$ws.Range("L69").Value = 18490.7139
$ws.Range("N69").Value = -20238.7139

$ws.Range("H72").Value = 3839.6592
$ws.Range("J72").Value = 6163.5713
$ws.Range("L72").Value = 55472.14169999999
$ws.Range("N72").Value = -64208.14169999999

$ws.Range("H116").Value = 1667.5
$ws.Range("I116").Value = 1563
$ws.Range("J116").Value = 1876.5
$ws.Range("K116").Value = 1563
$ws.Range("L116").Value = 1876.5
$ws.Range("M116").Value = 1879
$ws.Range("N116").Value = -8760.5

$ws.Range("H132").Value = 2268719.2
$ws.Range("I132").Value = 2422422.2
$ws.Range("J132").Value = 1599.5
$ws.Range("K132").Value = 7267266.600000001
$ws.Range("L132").Value = 4798.5
$ws.Range("M132").Value = -7264736.600000001
$ws.Range("N132").Value = -9858.5

$ws.Range("H138").Value = 4445927.5
$ws.Range("I138").Value = 7937206.5
$ws.Range("J138").Value = 2481.818
$ws.Range("K138").Value = 23811619.5
$ws.Range("L138").Value = 7445.454000000001
$ws.Range("M138").Value = -23806479.5
$ws.Range("N138").Value = -17725.454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13730.6045
$ws.Range("I32").Value = 13731.888
$ws.Range("K32").Value = 13731.888
$ws.Range("M32").Value = -13444.888

$ws.Range("H45").Value = 1217.6086
$ws.Range("I45").Value = 1006.8333
$ws.Range("K45").Value = 1006.8333
$ws.Range("M45").Value = -629.8333

$ws.Range("H61").Value = 1048
$ws.Range("I61").Value = 633.8125
$ws.Range("J61").Value = 3257
$ws.Range("K61").Value = 633.8125
$ws.Range("L61").Value = 3257
$ws.Range("M61").Value = -421.8125
$ws.Range("N61").Value = -3681

$ws.Range("H74").Value = 806.15
$ws.Range("I74").Value = 832.8823
$ws.Range("J74").Value = 654.6667
$ws.Range("K74").Value = 832.8823
$ws.Range("L74").Value = 654.6667
$ws.Range("M74").Value = 41.11770000000001
$ws.Range("N74").Value = -2402.6667

$ws.Range("H76").Value = 20000
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 20000
$ws.Range("N76").Value = -20676

$ws.Range("H77").Value = 806.15
$ws.Range("I77").Value = 832.8823
$ws.Range("J77").Value = 654.6667
$ws.Range("K77").Value = 4164.4115
$ws.Range("L77").Value = 3273.3335
$ws.Range("M77").Value = 203.5884999999998
$ws.Range("N77").Value = -12009.3335

$ws.Range("H79").Value = 20000
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 20000
$ws.Range("N79").Value = -22340

$ws.Range("H102").Value = 1782.7142
$ws.Range("I102").Value = 1624.5
$ws.Range("J102").Value = 1993.6666
$ws.Range("K102").Value = 1624.5
$ws.Range("L102").Value = 1993.6666
$ws.Range("M102").Value = -2.5
$ws.Range("N102").Value = -5237.6666

$ws.Range("H122").Value = 1722.2222
$ws.Range("I122").Value = 1722.2222
$ws.Range("K122").Value = 5166.6666
$ws.Range("M122").Value = -2716.6666

$ws.Range("H132").Value = 9570.111000000001
$ws.Range("I132").Value = 11439.615
$ws.Range("J132").Value = 4709.4
$ws.Range("K132").Value = 34318.845
$ws.Range("L132").Value = 14128.2
$ws.Range("M132").Value = -31788.845
$ws.Range("N132").Value = -19188.2

$ws.Range("H136").Value = 1048
$ws.Range("I136").Value = 633.8125
$ws.Range("J136").Value = 3257
$ws.Range("K136").Value = 1901.4375
$ws.Range("L136").Value = 9771
$ws.Range("M136").Value = 648.5625
$ws.Range("N136").Value = -14871

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 958.96295
$ws.Range("I58").Value = 767.8182
$ws.Range("K58").Value = 767.8182
$ws.Range("M58").Value = -564.8182

$ws.Range("H132").Value = 6023.375
$ws.Range("I132").Value = 6733.3335
$ws.Range("J132").Value = 5597.4
$ws.Range("K132").Value = 20200.0005
$ws.Range("L132").Value = 16792.2
$ws.Range("M132").Value = -17670.0005
$ws.Range("N132").Value = -21852.2

$ws.Range("H134").Value = 823.0526
$ws.Range("I134").Value = 747.56665
$ws.Range("J134").Value = 1106.125
$ws.Range("K134").Value = 2242.69995
$ws.Range("L134").Value = 3318.375
$ws.Range("M134").Value = 292.3000499999998
$ws.Range("N134").Value = -8388.375

$ws.Range("H136").Value = 958.96295
$ws.Range("I136").Value = 767.8182
$ws.Range("K136").Value = 2303.4546
$ws.Range("M136").Value = 246.5454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 342440.1
$ws.Range("I131").Value = 4964.5415
$ws.Range("J131").Value = 518514.28
$ws.Range("K131").Value = 14893.6245
$ws.Range("L131").Value = 1555542.84
$ws.Range("M131").Value = -9853.624500000002
$ws.Range("N131").Value = -1565622.84

$ws.Range("H132").Value = 1823.3914
$ws.Range("I132").Value = 1018.5833
$ws.Range("J132").Value = 2701.3635
$ws.Range("K132").Value = 9167.2497
$ws.Range("L132").Value = 24312.2715
$ws.Range("M132").Value = -6637.2497
$ws.Range("N132").Value = -29372.2715

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1800.091
$ws.Range("I126").Value = 2020.0769
$ws.Range("J126").Value = 1482.3334
$ws.Range("K126").Value = 6060.2307
$ws.Range("L126").Value = 4447.0002
$ws.Range("M126").Value = -3590.2307
$ws.Range("N126").Value = -9387.0002

$ws.Range("H132").Value = 23473.457
$ws.Range("I132").Value = 28452.432
$ws.Range("J132").Value = 3004.3333
$ws.Range("K132").Value = 85357.296
$ws.Range("L132").Value = 9012.999899999999
$ws.Range("M132").Value = -82827.296
$ws.Range("N132").Value = -14072.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1593.6666
$ws.Range("I16").Value = 998.1
$ws.Range("J16").Value = 2784.8
$ws.Range("K16").Value = 998.1
$ws.Range("L16").Value = 2784.8
$ws.Range("M16").Value = -828.1
$ws.Range("N16").Value = -3124.8

$ws.Range("H22").Value = 2210.3333
$ws.Range("I22").Value = 5500
$ws.Range("J22").Value = 565.5
$ws.Range("K22").Value = 5500
$ws.Range("L22").Value = 565.5
$ws.Range("M22").Value = -5205
$ws.Range("N22").Value = -1155.5

$ws.Range("H27").Value = 2210.3333
$ws.Range("I27").Value = 5500
$ws.Range("J27").Value = 565.5
$ws.Range("K27").Value = 5500
$ws.Range("L27").Value = 565.5
$ws.Range("M27").Value = -5393
$ws.Range("N27").Value = -779.5

$ws.Range("H122").Value = 2490.2903
$ws.Range("I122").Value = 2510.7368
$ws.Range("J122").Value = 2457.9167
$ws.Range("K122").Value = 7532.2104
$ws.Range("L122").Value = 7373.750100000001
$ws.Range("M122").Value = -5082.2104
$ws.Range("N122").Value = -12273.7501

$ws.Range("H136").Value = 4414.6665
$ws.Range("I136").Value = 6688.2354
$ws.Range("J136").Value = 1999
$ws.Range("K136").Value = 20064.7062
$ws.Range("L136").Value = 5997
$ws.Range("M136").Value = -17514.7062
$ws.Range("N136").Value = -11097

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 86408.414
$ws.Range("I122").Value = 103020.2
$ws.Range("J122").Value = 3349.5
$ws.Range("K122").Value = 309060.6
$ws.Range("L122").Value = 10048.5
$ws.Range("M122").Value = -306610.6
$ws.Range("N122").Value = -14948.5

$ws.Range("H126").Value = 9982
$ws.Range("I126").Value = 9982
$ws.Range("K126").Value = 29946
$ws.Range("M126").Value = -27476
